$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K5").Value = 384
$ws.Range("K9").Value = 547
$ws.Range("K10").Value = 415
$ws.Range("K12").Value = 645
$ws.Range("K14").Value = 393
$ws.Range("K15").Value = 561
$ws.Range("K20").Value = 387
$ws.Range("K22").Value = 554
$ws.Range("K26").Value = 366
$ws.Range("K34").Value = 552
$ws.Range("K36").Value = 730
$ws.Range("K39").Value = 563
$ws.Range("K42").Value = 355
